# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that any comma-separated part containing '@' (an email address) is moved
# to the front of the list, while preserving the relative order of the
# remaining parts (e.g. "System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    $withAt = @()
    $withoutAt = @()
    foreach ($part in $parts) {
        if ($part -like "*@*") {
            $withAt += $part
        } else {
            $withoutAt += $part
        }
    }

    $newParts = $withAt + $withoutAt
    $newVal = [string]::Join(", ", $newParts)

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
